$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map for this data refresh (prices, 1h volume %, and a few
# reshuffled coin rows - see commit message: symbol list update)
$updates = [ordered]@{
    "D2" = "306.40"
    "E2" = "-0.56%"
    "D3" = "38.90"
    "E3" = "7.01%"
    "D4" = "5.096"
    "E4" = "0.54%"
    "E5" = "-0.59%"
    "D6" = "1.949"
    "E6" = "-7.36%"
    "B7" = "KuCoinToken"
    "C7" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "D7" = "8.003"
    "E7" = "1.87%"
    "B8" = "MXToken"
    "C8" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8" = "0.9312"
    "E8" = "0.12%"
    "B9" = "LiechtensteinCryptoassetsExchange"
    "C9" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D9" = "0.1448"
    "E9" = "-0.40%"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1928"
    "E10" = "-0.18%"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D11" = "0.09092"
    "E11" = "0.24%"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "0.03502"
    "E12" = "1.34%"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D13" = "0.09797"
    "E13" = "-1.20%"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D14" = "0.001390"
    "E14" = "-1.39%"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D15" = "0.006024"
    "E15" = "-4.70%"
    "B16" = "LEO"
    "C16" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D16" = "3.785"
    "E16" = "-1.67%"
    "B17" = "GateToken"
    "C17" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D17" = "4.194"
    "E17" = "0.98%"
    "E18" = "2.04%"
    "E19" = "-1.29%"
    "E20" = "1.46%"
    "D21" = "4.792"
    "E21" = "-0.10%"
    "D22" = "0.2414"
    "E22" = "3.03%"
    "D23" = "0.04387"
    "E23" = "0.37%"
    "E24" = "0.52%"
    "D25" = "0.004276"
    "E25" = "-13.07%"
    "E26" = "0.13%"
    "D39" = "0.02032"
    "E39" = "0.89%"
    "D40" = "0.05042"
    "E40" = "-2.43%"
    "D41" = "0.007425"
    "E41" = "-0.81%"
    "D42" = "0.01020"
    "E42" = "0.50%"
    "D43" = "0.1350"
    "E43" = "-1.42%"
    "E44" = "-2.18%"
    "D45" = "0.009108"
    "E45" = "-8.55%"
    "D46" = "0.00006199"
    "E46" = "-1.28%"
    "E47" = "0.02%"
    "D48" = "0.003070"
    "E49" = "28.03%"
    "E50" = "0.02%"
    "E51" = "0.02%"
}

# Columns D (Price) and E (Volume %) store numeric-looking values as literal
# text, so force text format before assignment to stop Excel from parsing them
# into numbers (which would also strip formatting like trailing zeros / "%")
$textFormatCols = @("D", "E")

foreach ($ref in $updates.Keys) {
    $col = $ref -replace '[0-9]+$', ''
    if ($textFormatCols -contains $col) {
        $ws.Range($ref).NumberFormat = "@"
    }
    $ws.Range($ref).Value = $updates[$ref]
}
